# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45178 (2023-09-09) to 45179 (2023-09-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C369").Value = 45179
